# Apply "Add data for 2021-10-04" update:
# - Rename sheet/tab to reflect new "through" date (09-25 -> 09-26)
# - Update the September row label text accordingly
# - Update September row values (row 10) and Total row values (row 11)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (tab) - mirrors the <sheet name="..."> change
$ws.Name = "Through 2021-09-26"

# Update the row label for September in column A (row 10)
$ws.Range("A10").Value = "September (through 09-26)"

# Update September row (row 10) values for years 2015-2021 (columns B-H)
$ws.Range("B10").Value = 27
$ws.Range("C10").Value = 39
$ws.Range("D10").Value = 67
$ws.Range("E10").Value = 50
$ws.Range("F10").Value = 63
$ws.Range("G10").Value = 99
$ws.Range("H10").Value = 164

# Update Total row (row 11) values for years 2015-2021 (columns B-H)
$ws.Range("B11").Value = 221
$ws.Range("C11").Value = 420
$ws.Range("D11").Value = 618
$ws.Range("E11").Value = 540
$ws.Range("F11").Value = 412
$ws.Range("G11").Value = 883
$ws.Range("H11").Value = 1234
